$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Feria Lagunitas de Puerto Montt -
# Brocoli. It belongs (chronologically) right before the current row 207, so
# insert a fresh row there; Excel shifts the existing rows 207:253 down to
# 208:254 and carries the date-column (D) number format along automatically.
$ws.Rows.Item(207).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A207").Value = 4
$ws.Range("B207").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C207").Value = "Los Lagos"
$ws.Range("D207").Value = 44543
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = 100112023
$ws.Range("G207").Value = "Brócoli"
$ws.Range("H207").Value = "Sin especificar"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 500
$ws.Range("K207").Value = 1200
$ws.Range("L207").Value = 1200
$ws.Range("M207").Value = 1200
$ws.Range("N207").Value = "$/unidad"
$ws.Range("O207").Value = "Región Metropolitana"
$ws.Range("P207").Value = 1200
$ws.Range("Q207").Value = 1
$ws.Range("R207").Value = "Hortaliza"
